$d = $word.ActiveDocument

# 1. Merge the two runs of the first list item into a single run.
$findRange = $d.Content
[void]$findRange.Find.Execute(
    "Writing before a log file has been specified. It should return a custom error object indicating the issue.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Writing before a log file has been specified. It should return a custom error object indicating the issue.",
    2)

# 2. Remove the _GoBack bookmark from that paragraph (it gets re-added later).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3. Append a new list paragraph describing the GetSize test case.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara.Range.InsertParagraphAfter()
$getSizePara = $d.Paragraphs.Item($d.Paragraphs.Count)
$getSizeXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>GetSize</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> on a file that does not exist ' + [char]0x2013 + ' custom error</w:t></w:r>' +
    '</w:p>'
[void]$getSizePara.Range.InsertXML($getSizeXml)

# 4. Append a final empty list paragraph that now hosts the _GoBack bookmark.
$lastPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara2.Range.InsertParagraphAfter()
$bmPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
[void]$bmPara.Range.InsertXML($bmXml)
